$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the previous full extent (A1:AD19) ---
# The old sheet had 30 columns (B:AD) with duplicate blocks in K:T and U:AD;
# the new layout only uses columns A:T, so we wipe the old range completely
# (content + formatting) before writing the new grid.
$ws.Range("A1:AD30").Clear()

# --- Target grid: row number -> { column index -> value } ---
$grid = @{
    1 = @{ 2=0; 3=1; 4=2; 5=3; 6=4; 7=5; 8=6; 9=7; 10=8; 11=9; 12=10; 13=11; 14=12; 15=13; 16=14; 17=15; 18=16; 19=17; 20=18 }
    2 = @{ 1=0; 2='HKL'; 3='[2, 2, 0]'; 4='[2, 0, 0]'; 5='[2, 1, 1]'; 6='[4, 0, 0]'; 7='[3, 2, 1]'; 8='[1, 1, 0]'; 9='[2, 2, 2]'; 10='[3, 1, 0]'; 11='1Pair-A'; 12='1Pair-B'; 13='2Pairs-A'; 14='2Pairs-B'; 15='3Pairs-A'; 16='3Pairs-B'; 17='3Pairs-C'; 18='4Pairs'; 19='5A4F'; 20='MaxUnique' }
    3 = @{ 1=1; 2='BT8Hex_2.5'; 3=0.9984165034208764; 4=1.006682875984724; 5=0.998280824243695; 6=1.006682875984724; 7=0.9983203792005639; 8=0.9984165034208764; 9=0.9953896846980712; 10=1.003706984602589; 11=0.9984165034208764; 12=0.998280824243695; 13=1.002481850114209; 14=1.002481850114209; 15=1.002890228277002; 16=1.001126734549765; 17=1.001126734549765; 18=1.000449176767543; 19=1.000449176767543; 20=1.00013287535842 }
    4 = @{ 1=2; 2='BT8Hex_5'; 3=0.9969422320005472; 4=1.012901734907552; 5=0.996681417728606; 6=1.012901734907552; 7=0.9967574553995612; 8=0.9969422320005472; 9=0.9911007671395456; 10=1.007156316701217; 11=0.9969422320005472; 12=0.996681417728606; 13=1.004791576318079; 14=1.004791576318079; 15=1.005579823112458; 16=1.002175128212235; 17=1.002175128212235; 18=1.000866904159313; 19=1.000866904159313; 20=1.000256653979505 }
    5 = @{ 1=3; 2='BT8Hex_10'; 3=0.994123307797587; 4=1.024734959273292; 5=0.9936456756579407; 6=1.024734959273292; 7=0.9937849243453009; 8=0.994123307797587; 9=0.982964154595072; 10=1.013714773337744; 11=0.994123307797587; 12=0.9936456756579407; 13=1.009190317465616; 14=1.009190317465616; 15=1.010698469422992; 16=1.004167980909606; 17=1.004167980909606; 18=1.001656812631601; 19=1.001656812631601; 20=1.000494632501156 }
    6 = @{ 1=4; 2='BT8Hex_15'; 3=0.991371668182053; 4=1.036283888575724; 5=0.9906831015495905; 6=1.036283888575724; 7=0.9908838550960989; 8=0.991371668182053; 9=0.9750237956878859; 10=1.020115499873508; 11=0.991371668182053; 12=0.9906831015495905; 13=1.013483495062657; 14=1.013483495062657; 15=1.015694163332941; 16=1.006112886102456; 17=1.006112886102456; 18=1.002427581622355; 19=1.002427581622355; 20=1.00072696816081 }
    7 = @{ 1=5; 2='Spiral2.5'; 3=0.9998827501451563; 4=1.000446266389896; 5=0.9998915905350169; 6=1.000446266389896; 7=0.9998890118568226; 8=0.9998827501451563; 9=0.9997125897489934; 10=1.000243402923784; 11=0.9998827501451563; 12=0.9998915905350169; 13=1.000168928462456; 14=1.000168928462456; 15=1.000193753282899; 16=1.000073535690023; 17=1.000073535690023; 18=1.000025839303806; 19=1.000025839303806; 20=1.000010935266612 }
    8 = @{ 1=6; 2='Spiral5'; 3=0.9997300033019333; 4=1.001037090759158; 5=0.9997466848079257; 6=1.001037090759158; 7=0.9997418200285424; 8=0.9997300033019333; 9=0.999327667949438; 10=1.000566541628783; 11=0.9997300033019333; 12=0.9997466848079257; 13=1.000391887783542; 14=1.000391887783542; 15=1.000450105731956; 16=1.000171259623006; 17=1.000171259623006; 18=1.000060945542738; 19=1.000060945542738; 20=1.000024968079297 }
    9 = @{ 1=7; 2='Spiral7.5'; 3=0.9996493323894972; 4=1.001375884524083; 5=0.9996597484487062; 6=1.001375884524083; 7=0.9996567106590847; 8=0.9996493323894972; 9=0.9990946448547211; 10=1.000754328260677; 11=0.9996493323894972; 12=0.9996597484487062; 13=1.000517816486394; 14=1.000517816486394; 15=1.000596653744489; 16=1.000228321787429; 17=1.000228321787429; 18=1.000083574437946; 19=1.000083574437946; 20=1.000031774856128 }
    10 = @{ 1=8; 2='Spiral10'; 3=0.9992142521063629; 4=1.003043067328038; 5=0.9992531112813167; 6=1.003043067328038; 7=0.9992417806078145; 8=0.9992142521063629; 9=0.998015698339994; 10=1.001664696742958; 11=0.9992142521063629; 12=0.9992531112813167; 13=1.001148089304678; 14=1.001148089304678; 15=1.001320291784104; 16=1.000503476905239; 17=1.000503476905239; 18=1.00018117070552; 19=1.00018117070552; 20=1.000072101067747 }
    11 = @{ 1=9; 2='Spiral15'; 3=0.9987764178989033; 4=1.004827798530688; 5=0.9988022880813106; 6=1.004827798530688; 7=0.9987947453535956; 8=0.9987764178989033; 9=0.9968110297997075; 10=1.002649307173324; 11=0.9987764178989033; 12=0.9988022880813106; 13=1.001815043305999; 14=1.001815043305999; 15=1.002093131261774; 16=1.000802168170301; 17=1.000802168170301; 18=1.000295730602451; 19=1.000295730602451; 20=1.000110264472922 }
    12 = @{ 1=10; 2='OffsetF45'; 3=1.014369839509032; 4=0.9637354397832227; 5=1.006119641121082; 6=0.9637354397832227; 7=1.008524948117354; 8=1.014369839509032; 9=1.014747592515798; 10=0.9819638286793322; 11=1.014369839509032; 12=1.006119641121082; 13=0.9849275404521525; 14=0.9849275404521525; 15=0.9839396365278791; 16=0.994741640137779; 17=0.994741640137779; 18=0.9996486899805923; 19=0.9996486899805923; 20=0.9982435482876367 }
    13 = @{ 1=11; 2='OffsetA45'; 3=0.9994418808283211; 4=0.993306383595565; 5=1.002913139413011; 6=0.993306383595565; 7=1.001901106146708; 8=0.9994418808283211; 9=1.008429562082902; 10=0.9955151642948434; 11=0.9994418808283211; 12=1.002913139413011; 13=0.9981097615042878; 14=0.9981097615042878; 15=0.9972448957678063; 16=0.9985538012789655; 17=0.9985538012789655; 18=0.9987758211663044; 19=0.9987758211663044; 20=1.000251206060225 }
    14 = @{ 1=12; 2='OffsetFTD'; 3=1.016711760342585; 4=0.92505585301876; 5=1.019860629219225; 6=0.92505585301876; 7=1.018942583935212; 8=1.016711760342585; 9=1.053561470836624; 10=0.9580519823176037; 11=1.016711760342585; 12=1.019860629219225; 13=0.9724582411189927; 14=0.9724582411189927; 15=0.9676561548518631; 16=0.9872094141935234; 17=0.9872094141935234; 18=0.9945850007307886; 19=0.9945850007307886; 20=0.9986973799450016 }
    15 = @{ 1=13; 2='OffsetATD'; 3=1.005927421620521; 4=0.9796151116138746; 5=1.004634448438814; 6=0.9796151116138746; 7=1.00501141069853; 8=1.005927421620521; 9=1.012112248171304; 10=0.989087548326328; 11=1.005927421620521; 12=1.004634448438814; 13=0.9921247800263444; 14=0.9921247800263444; 15=0.9911123694596723; 16=0.9967256605577366; 17=0.9967256605577366; 18=0.9990261008234327; 19=0.9990261008234327; 20=0.9993980314782287 }
    16 = @{ 1=14; 2='Holden2.5'; 3=0.9688758124299928; 4=1.130475756588137; 5=0.9665506452570605; 6=1.130475756588137; 7=0.9672285296489703; 8=0.9688758124299928; 9=0.9103588222141746; 10=1.072299786732198; 11=0.9688758124299928; 12=0.9665506452570605; 13=1.048513200922599; 14=1.048513200922599; 15=1.056442062859132; 16=1.021967404758397; 17=1.021967404758397; 18=1.008694506676296; 19=1.008694506676296; 20=1.002631558811756 }
    17 = @{ 1=15; 2='Holden5'; 3=0.974010769161123; 4=1.107021749341802; 5=0.9728188935891064; 6=1.107021749341802; 7=0.9731663740820091; 8=0.974010769161123; 9=0.9272900186076858; 10=1.059137804364969; 11=0.974010769161123; 12=0.9728188935891064; 13=1.039920321465454; 14=1.039920321465454; 15=1.046326149098626; 16=1.017950470697344; 17=1.017950470697344; 18=1.006965545313288; 19=1.006965545313288; 20=1.002240934857783 }
    18 = @{ 1=16; 2='Holden10'; 3=0.9844069256841594; 4=1.059691321762264; 5=0.9854495268943106; 6=1.059691321762264; 7=0.9851455617736108; 8=0.9844069256841594; 9=0.9613973261401062; 10=1.032588945575998; 11=0.9844069256841594; 12=0.9854495268943106; 13=1.022570424328287; 14=1.022570424328287; 15=1.025909931410858; 16=1.009849258113578; 17=1.009849258113578; 18=1.003488675006223; 19=1.003488675006223; 20=1.001446601305075 }
    19 = @{ 1=17; 2='Holden15'; 3=0.9826808153397703; 4=1.065273601121967; 5=0.9842375901251185; 6=1.065273601121967; 7=0.9837837174437928; 8=0.9826808153397703; 9=0.958263434481346; 10=1.035540200161451; 11=0.9826808153397703; 12=0.9842375901251185; 13=1.024755595623543; 14=1.024755595623543; 15=1.028350463802846; 16=1.010730668862285; 17=1.010730668862285; 18=1.003718205481657; 19=1.003718205481657; 20=1.001629893112241 }
    20 = @{ 1=18; 2='HexGrid-90degTilt2.5degRes'; 3=0.9999837951951559; 4=1.000033542013678; 5=0.9999959583051831; 6=1.000033542013678; 7=0.9999924109408161; 8=0.9999837951951559; 9=0.999991536818643; 10=1.000015635461891; 11=0.9999837951951559; 12=0.9999959583051831; 13=1.00001475015943; 14=1.00001475015943; 15=1.000015045260251; 16=1.000004431838005; 17=1.000004431838005; 18=0.999999272677293; 19=0.999999272677293; 20=1.000002146455895 }
    21 = @{ 1=19; 2='HexGrid-90degTilt5degRes'; 3=0.9998503731134427; 4=1.000591287493226; 5=0.999853180550645; 6=1.000591287493226; 7=0.9998523610804048; 8=0.9998503731134427; 9=0.9996090144214869; 10=1.000324561050013; 11=0.9998503731134427; 12=0.999853180550645; 13=1.000222234021936; 14=1.000222234021936; 15=1.000256343031295; 16=1.000098280385771; 17=1.000098280385771; 18=1.000036303567689; 19=1.000036303567689; 20=1.000013462951536 }
    22 = @{ 1=20; 2='HexGrid-90degTilt10degRes'; 3=0.9995017206906265; 4=1.002064203086009; 5=0.9994740786236306; 6=1.002064203086009; 7=0.9994821330995862; 8=0.9995017206906265; 9=0.9985922663097387; 10=1.001141712634882; 11=0.9995017206906265; 12=0.9994740786236306; 13=1.00076914085482; 14=1.00076914085482; 15=1.000893331448174; 16=1.000346667466755; 17=1.000346667466755; 18=1.000135430772723; 19=1.000135430772723; 20=1.000042685740745 }
    23 = @{ 1=21; 2='HexGrid-90degTilt15degRes'; 3=0.998862582397678; 4=1.004817411921165; 5=0.9987584562905129; 6=1.004817411921165; 7=0.9987888140835085; 8=0.998862582397678; 9=0.9966693877618064; 10=1.002673680877842; 11=0.998862582397678; 12=0.9987584562905129; 13=1.001787934105839; 14=1.001787934105839; 15=1.00208318302984; 16=1.000812816869785; 17=1.000812816869785; 18=1.000325258251759; 19=1.000325258251759; 20=1.000095055555419 }
}

foreach ($rowNum in $grid.Keys) {
    $rowCells = $grid[$rowNum]
    foreach ($colNum in $rowCells.Keys) {
        $ws.Cells.Item([int]$rowNum, [int]$colNum).Value = $rowCells[$colNum]
    }
}

# --- Reapply header/index formatting: bold, centered, thin border ---
# (Border set first, then font/alignment, to match the single pre-existing
#  style record instead of generating extra unused style entries.)
$headerRange = $ws.Range("B1:T1")
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexRange = $ws.Range("A2:A23")
$indexRange.Borders.LineStyle = 1
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
